$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 16160672
$ws.Range("I19").Value = 15653096
$ws.Range("J19").Value = 16668249
$ws.Range("K19").Value = 15653096
$ws.Range("L19").Value = 16668249
$ws.Range("M19").Value = -15652921
$ws.Range("N19").Value = -16668599

$ws.Range("H33").Value = 300
$ws.Range("I33").Value = 116.666664
$ws.Range("K33").Value = 116.666664
$ws.Range("M33").Value = 112.333336

$ws.Range("H43").Value = 902.9286
$ws.Range("I43").Value = 675.25
$ws.Range("J43").Value = 994
$ws.Range("K43").Value = 675.25
$ws.Range("L43").Value = 994
$ws.Range("M43").Value = -606.25
$ws.Range("N43").Value = -1132

$ws.Range("H76").Value = 3577.9443
$ws.Range("I76").Value = 3500.5386
$ws.Range("J76").Value = 3779.2
$ws.Range("K76").Value = 3500.5386
$ws.Range("L76").Value = 3779.2
$ws.Range("M76").Value = -3185.5386
$ws.Range("N76").Value = -4409.2

$ws.Range("H79").Value = 3577.9443
$ws.Range("I79").Value = 3500.5386
$ws.Range("J79").Value = 3779.2
$ws.Range("K79").Value = 3500.5386
$ws.Range("L79").Value = 3779.2
$ws.Range("M79").Value = -2408.5386
$ws.Range("N79").Value = -5963.2

$ws.Range("H116").Value = 2201808
$ws.Range("I116").Value = 5921359.5
$ws.Range("J116").Value = 3891.182
$ws.Range("K116").Value = 5921359.5
$ws.Range("L116").Value = 3891.182
$ws.Range("M116").Value = -5917917.5
$ws.Range("N116").Value = -10775.182

$ws.Range("H125").Value = 975.5714
$ws.Range("I125").Value = 599.8
$ws.Range("J125").Value = 1184.3334
$ws.Range("K125").Value = 5398.2
$ws.Range("L125").Value = 10659.0006
$ws.Range("M125").Value = -2938.2
$ws.Range("N125").Value = -15579.0006

$ws.Range("H132").Value = 1918.8478
$ws.Range("I132").Value = 1895.9722
$ws.Range("J132").Value = 2001.2
$ws.Range("K132").Value = 5687.9166
$ws.Range("L132").Value = 6003.6
$ws.Range("M132").Value = -3157.9166
$ws.Range("N132").Value = -11063.6

$ws.Range("H137").Value = 1992.5264
$ws.Range("I137").Value = 930
$ws.Range("K137").Value = 2790
$ws.Range("M137").Value = -240

$ws.Range("H138").Value = 3681.6042
$ws.Range("I138").Value = 1109.9584
$ws.Range("J138").Value = 4538.8193
$ws.Range("K138").Value = 3329.8752
$ws.Range("L138").Value = 13616.4579
$ws.Range("M138").Value = 1810.1248
$ws.Range("N138").Value = -23896.4579

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24396288
$ws.Range("I32").Value = 28574840
$ws.Range("K32").Value = 28574840
$ws.Range("M32").Value = -28574553

$ws.Range("H45").Value = 1682.826
$ws.Range("I45").Value = 934
$ws.Range("J45").Value = 1795.15
$ws.Range("K45").Value = 934
$ws.Range("L45").Value = 1795.15
$ws.Range("M45").Value = -557
$ws.Range("N45").Value = -2549.15

$ws.Range("H74").Value = 3480.1282
$ws.Range("I74").Value = 4134.2
$ws.Range("J74").Value = 1299.8889
$ws.Range("K74").Value = 4134.2
$ws.Range("L74").Value = 1299.8889
$ws.Range("M74").Value = -3260.2
$ws.Range("N74").Value = -3047.8889

$ws.Range("H77").Value = 3480.1282
$ws.Range("I77").Value = 4134.2
$ws.Range("J77").Value = 1299.8889
$ws.Range("K77").Value = 20671
$ws.Range("L77").Value = 6499.4445
$ws.Range("M77").Value = -16303
$ws.Range("N77").Value = -15235.4445

$ws.Range("H88").Value = 2000
$ws.Range("I88").Value = 2000
$ws.Range("K88").Value = 2000
$ws.Range("M88").Value = -1594

$ws.Range("H91").Value = 2000
$ws.Range("I91").Value = 2000
$ws.Range("K91").Value = 2000
$ws.Range("M91").Value = -596

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1137.0454
$ws.Range("I94").Value = 1225.9375
$ws.Range("J94").Value = 900
$ws.Range("K94").Value = 1225.9375
$ws.Range("L94").Value = 900
$ws.Range("M94").Value = -774.9375
$ws.Range("N94").Value = -1802

$ws.Range("H134").Value = 2308.2
$ws.Range("I134").Value = 2161.7646
$ws.Range("J134").Value = 3138
$ws.Range("K134").Value = 6485.293799999999
$ws.Range("L134").Value = 9414
$ws.Range("M134").Value = -3950.293799999999
$ws.Range("N134").Value = -14484

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 89
$ws.Range("I7").Value = 88.2
$ws.Range("J7").Value = 90
$ws.Range("K7").Value = 88.2
$ws.Range("L7").Value = 90
$ws.Range("M7").Value = 24.8
$ws.Range("N7").Value = -316

$ws.Range("H31").Value = 2032.6285
$ws.Range("I31").Value = 1659.1034
$ws.Range("J31").Value = 3838
$ws.Range("K31").Value = 1659.1034
$ws.Range("L31").Value = 3838
$ws.Range("M31").Value = -1364.1034
$ws.Range("N31").Value = -4428

$ws.Range("H34").Value = 2032.6285
$ws.Range("I34").Value = 1659.1034
$ws.Range("J34").Value = 3838
$ws.Range("K34").Value = 1659.1034
$ws.Range("L34").Value = 3838
$ws.Range("M34").Value = -1457.1034
$ws.Range("N34").Value = -4242

$ws.Range("H58").Value = 1086.7
$ws.Range("I58").Value = 892.619
$ws.Range("J58").Value = 1539.5555
$ws.Range("K58").Value = 892.619
$ws.Range("L58").Value = 1539.5555
$ws.Range("M58").Value = -689.619
$ws.Range("N58").Value = -1945.5555

$ws.Range("H86").Value = 6574.625
$ws.Range("I86").Value = 8949.75
$ws.Range("J86").Value = 4199.5
$ws.Range("K86").Value = 8949.75
$ws.Range("L86").Value = 4199.5
$ws.Range("M86").Value = -7826.75
$ws.Range("N86").Value = -6445.5

$ws.Range("H89").Value = 6574.625
$ws.Range("I89").Value = 8949.75
$ws.Range("J89").Value = 4199.5
$ws.Range("K89").Value = 44748.75
$ws.Range("L89").Value = 20997.5
$ws.Range("M89").Value = -39132.75
$ws.Range("N89").Value = -32229.5

$ws.Range("H95").Value = 38330
$ws.Range("J95").Value = 38330
$ws.Range("L95").Value = 38330
$ws.Range("N95").Value = -43822

$ws.Range("H132").Value = 3038.2258
$ws.Range("I132").Value = 2170.6667
$ws.Range("J132").Value = 3851.5625
$ws.Range("K132").Value = 6512.000100000001
$ws.Range("L132").Value = 11554.6875
$ws.Range("M132").Value = -3982.000100000001
$ws.Range("N132").Value = -16614.6875

$ws.Range("H134").Value = 2251.4546
$ws.Range("I134").Value = 1551.25
$ws.Range("J134").Value = 3328.6924
$ws.Range("K134").Value = 4653.75
$ws.Range("L134").Value = 9986.0772
$ws.Range("M134").Value = -2118.75
$ws.Range("N134").Value = -15056.0772

$ws.Range("H136").Value = 1086.7
$ws.Range("I136").Value = 892.619
$ws.Range("J136").Value = 1539.5555
$ws.Range("K136").Value = 2677.857
$ws.Range("L136").Value = 4618.666499999999
$ws.Range("M136").Value = -127.857
$ws.Range("N136").Value = -9718.666499999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5114.391
$ws.Range("I131").Value = 526.8570999999999
$ws.Range("J131").Value = 7121.4375
$ws.Range("K131").Value = 1580.5713
$ws.Range("L131").Value = 21364.3125
$ws.Range("M131").Value = 3459.4287
$ws.Range("N131").Value = -31444.3125

$ws.Range("H132").Value = 2526529.5
$ws.Range("I132").Value = 1304
$ws.Range("K132").Value = 11736
$ws.Range("M132").Value = -9206

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 6610.421
$ws.Range("I113").Value = 885.7778
$ws.Range("K113").Value = 885.7778
$ws.Range("M113").Value = 1284.2222

$ws.Range("H126").Value = 3911.625
$ws.Range("I126").Value = 3472.3076
$ws.Range("J126").Value = 4430.8184
$ws.Range("K126").Value = 10416.9228
$ws.Range("L126").Value = 13292.4552
$ws.Range("M126").Value = -7946.9228
$ws.Range("N126").Value = -18232.4552

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 407
$ws.Range("J22").Value = 478.33334
$ws.Range("L22").Value = 478.33334
$ws.Range("N22").Value = -1068.33334

$ws.Range("H27").Value = 407
$ws.Range("J27").Value = 478.33334
$ws.Range("L27").Value = 478.33334
$ws.Range("N27").Value = -692.33334

$ws.Range("H100").Value = 589899.5600000001
$ws.Range("I100").Value = 3334763.2
$ws.Range("J100").Value = 1714.5
$ws.Range("K100").Value = 3334763.2
$ws.Range("L100").Value = 1714.5
$ws.Range("M100").Value = -3334222.2
$ws.Range("N100").Value = -2796.5

$ws.Range("H133").Value = 52122.5
$ws.Range("J133").Value = 52122.5
$ws.Range("L133").Value = 52122.5
$ws.Range("N133").Value = -57182.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 5548.84
$ws.Range("I136").Value = 5610.8
$ws.Range("J136").Value = 5301
$ws.Range("K136").Value = 16832.4
$ws.Range("L136").Value = 15903
$ws.Range("M136").Value = -14282.4
$ws.Range("N136").Value = -21003
